# Apply updated cryptocurrency market data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D on the rows below hold plain decimal-looking text (e.g. "0.999").
# Force them to keep their original Text formatting so Excel does not
# silently convert them into numeric values (which would also mangle
# values such as "0.0000174" into scientific notation).
$textCells = @("D4", "D5", "D6", "D8", "D10", "D12", "D14", "D16", "D19", "D20", "D21", "D22", "D23", "D25", "D27", "D28", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D39", "D40", "D42", "D43", "D45", "D46", "D47", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "66.845.41"
$ws.Range("E2").Value = "  -1.97%  "
$ws.Range("D3").Value = "2.482.73"
$ws.Range("E3").Value = "  -2.06%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "585.45"
$ws.Range("E5").Value = "  -1.45%  "
$ws.Range("D6").Value = "168.67"
$ws.Range("E6").Value = "  -4.61%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "0.516"
$ws.Range("E8").Value = "  -3.02%  "
$ws.Range("D9").Value = "2.480.73"
$ws.Range("E9").Value = "  -2.09%  "
$ws.Range("D10").Value = "0.136"
$ws.Range("E10").Value = "  -4.52%  "
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("D12").Value = "4.95"
$ws.Range("E12").Value = "  -4.00%  "
$ws.Range("E13").Value = "  -3.17%  "
$ws.Range("D14").Value = "25.93"
$ws.Range("E14").Value = "  -3.96%  "
$ws.Range("D15").Value = "2.930.63"
$ws.Range("E15").Value = "  -2.14%  "
$ws.Range("D16").Value = "0.0000174"
$ws.Range("E16").Value = "  -3.20%  "
$ws.Range("D17").Value = "66.671.33"
$ws.Range("E17").Value = "  -1.82%  "
$ws.Range("D18").Value = "2.479.28"
$ws.Range("E18").Value = "  -1.43%  "
$ws.Range("D19").Value = "11.68"
$ws.Range("E19").Value = "  +1.52%  "
$ws.Range("D20").Value = "7.93"
$ws.Range("E20").Value = "  -1.34%  "
$ws.Range("D21").Value = "363.33"
$ws.Range("E21").Value = "  -0.64%  "
$ws.Range("D22").Value = "4.06"
$ws.Range("E22").Value = "  -3.88%  "
$ws.Range("D23").Value = "4.46"
$ws.Range("E23").Value = "  -5.30%  "
$ws.Range("E24").Value = "  +0.21%  "
$ws.Range("D25").Value = "70.85"
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("E26").Value = "  -6.25%  "
$ws.Range("D27").Value = "9.46"
$ws.Range("E27").Value = "  -8.04%  "
$ws.Range("D28").Value = "0.997"
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("D29").Value = "2.603.18"
$ws.Range("E29").Value = "  -2.46%  "
$ws.Range("D30").Value = "0.0₃0932"
$ws.Range("E30").Value = "  -6.33%  "
$ws.Range("D31").Value = "8.09"
$ws.Range("E31").Value = "  -2.40%  "
$ws.Range("D32").Value = "517.39"
$ws.Range("E32").Value = "  -5.99%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "1.28"
$ws.Range("E33").Value = "  -5.49%  "
$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").Value = "1.83"
$ws.Range("E34").Value = "  -2.71%  "
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").Value = "0.127"
$ws.Range("E36").Value = "  -2.54%  "
$ws.Range("D37").Value = "158.21"
$ws.Range("E37").Value = "  +1.02%  "
$ws.Range("E38").Value = "  -3.81%  "
$ws.Range("D39").Value = "18.96"
$ws.Range("E39").Value = "  +0.29%  "
$ws.Range("D40").Value = "18.55"
$ws.Range("E40").Value = "  -0.65%  "
$ws.Range("E41").Value = "  -3.58%  "
$ws.Range("D42").Value = "4.96"
$ws.Range("E42").Value = "  -4.93%  "
$ws.Range("D43").Value = "0.334"
$ws.Range("E43").Value = "  -6.72%  "
$ws.Range("E44").Value = "  -3.01%  "
$ws.Range("D45").Value = "39.24"
$ws.Range("E45").Value = "  -2.07%  "
$ws.Range("D46").Value = "143.14"
$ws.Range("E46").Value = "  -3.09%  "
$ws.Range("D47").Value = "0.538"
$ws.Range("E47").Value = "  -4.79%  "
$ws.Range("E48").Value = "  -3.94%  "
$ws.Range("D49").Value = "0.0₆0269"
$ws.Range("E49").Value = "  -3.64%  "
$ws.Range("E50").Value = "  -3.32%  "
$ws.Range("D51").Value = "0.0738"
$ws.Range("E51").Value = "  -2.75%  "
